$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.775.71"
$ws.Range("E2").Value = "  +1.98%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.859.44"
$ws.Range("E3").Value = "  +1.64%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.47"
$ws.Range("E5").Value = "  +0.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6436"
$ws.Range("E6").Value = "  +4.35%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.12"
$ws.Range("E8").Value = "  +5.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07527"
$ws.Range("E9").Value = "  +2.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.2985"
$ws.Range("E10").Value = "  +1.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "24.63"
$ws.Range("E11").Value = "  +6.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07694"
$ws.Range("E12").Value = "  +0.54%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.860.15"
$ws.Range("E13").Value = "  +1.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.052"
$ws.Range("E14").Value = "  +1.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6920"
$ws.Range("E15").Value = "  +2.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "84.03"
$ws.Range("E16").Value = "  +1.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009910"
$ws.Range("E17").Value = "  +10.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.136"
$ws.Range("E18").Value = "  +4.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.797.30"
$ws.Range("E19").Value = "  +2.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.113.23"
$ws.Range("E20").Value = "  +1.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "236.68"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.65"
$ws.Range("E22").Value = "  +1.19%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.564"
$ws.Range("E24").Value = "  +2.33%  "

$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.57"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1424"
$ws.Range("E27").Value = "  +2.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.569"
$ws.Range("E28").Value = "  +0.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.93"
$ws.Range("E29").Value = "  +1.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06232"
$ws.Range("E30").Value = "  +7.67%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.495"
$ws.Range("E31").Value = "  +0.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.290"
$ws.Range("E32").Value = "  +5.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.157"
$ws.Range("E33").Value = "  +1.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.098"
$ws.Range("E34").Value = "  +0.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.900"
$ws.Range("E35").Value = "  +2.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.174"
$ws.Range("E36").Value = "  +3.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7324"
$ws.Range("E37").Value = "  +1.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.609"
$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("E39").Value = "  -1.46%  "

$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.222.87"
$ws.Range("E40").Value = "  +0.17%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01790"
$ws.Range("E41").Value = "  +1.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.311"
$ws.Range("E42").Value = "  +1.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9201"
$ws.Range("E43").Value = "  +1.52%  "

$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.023.15"
$ws.Range("E45").Value = "  +0.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.07"
$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "67.13"
$ws.Range("E47").Value = "  +2.37%  "

$ws.Range("E48").Value = "  +0.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4069"
$ws.Range("E49").Value = "  +1.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.197"
$ws.Range("E50").Value = "  +0.55%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.671"
$ws.Range("E51").Value = "  +5.46%  "
